# Actualizacion automatica: se corrigen los valores de ventas de diciembre
# en las hojas "VENTAS POR GRUPO" y "VENTA MENSUAL".

$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Fila 4 (AGUILAR REYES CESAR VINICIO)
$wsGrupo.Range("D4").Value = 183.17
$wsGrupo.Range("E4").Value = 95.48999999999999
$wsGrupo.Range("M4").Value = 1360.13
$wsGrupo.Range("N4").Value = 234.38

# Fila 24 (FEIJOO MARIN MAURICIO ENRIQUE)
$wsGrupo.Range("M24").Value = 3550.13

# Fila 36 (ORTEGA ROMAN KLEBER ERWIN)
$wsGrupo.Range("M36").Value = 12682.51
$wsGrupo.Range("O36").Value = 1058.37

# Fila 37 (ORTEGA ROMAN LUIS FERNANDO)
$wsGrupo.Range("C37").Value = 388.8
$wsGrupo.Range("L37").Value = 506.88
$wsGrupo.Range("M37").Value = 3409.47

# Fila 48 (RUIZ TINIZARAY YOHANNA MARYURI)
$wsGrupo.Range("D48").Value = 475.2

# Fila 52 (TORO BLACIO MARIA DEL CISNE)
$wsGrupo.Range("D52").Value = 475.2
$wsGrupo.Range("O52").Value = 1058.37

# Fila 56 (conteo "X de 54" por columna, recalculado tras los nuevos valores)
$wsGrupo.Range("C56").Value = "3 de 54"
$wsGrupo.Range("D56").Value = "6 de 54"
$wsGrupo.Range("E56").Value = "2 de 54"
$wsGrupo.Range("L56").Value = "4 de 54"
$wsGrupo.Range("N56").Value = "1 de 54"
$wsGrupo.Range("O56").Value = "2 de 54"

# --- Hoja "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

$wsMensual.Range("F4").Value = 2132.37
$wsMensual.Range("F24").Value = 3773
$wsMensual.Range("F36").Value = 14866.21
$wsMensual.Range("F37").Value = 6473.25
$wsMensual.Range("F48").Value = 745.2
$wsMensual.Range("F53").Value = 1533.57
$wsMensual.Range("F54").Value = 1533.57
$wsMensual.Range("F60").Value = 50284.32
